$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 548
$ws.Cells.Item(3, 6).Value = 10348
$ws.Cells.Item(4, 6).Value = 227
$ws.Cells.Item(5, 6).Value = 94
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 133
$ws.Cells.Item(9, 6).Value = 12042
$ws.Cells.Item(10, 6).Value = 12615
$ws.Cells.Item(11, 6).Value = 1299
$ws.Cells.Item(12, 6).Value = 1275
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(20, 6).Value = 1991
$ws.Cells.Item(21, 6).Value = 1012
$ws.Cells.Item(23, 6).Value = 878
$ws.Cells.Item(24, 6).Value = 13
$ws.Cells.Item(26, 6).Value = 714
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 2003
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(31, 6).Value = 1669
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 33
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 3660
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 124
$ws.Cells.Item(43, 6).Value = 35
$ws.Cells.Item(45, 6).Value = 281
$ws.Cells.Item(46, 6).Value = 18
$ws.Cells.Item(48, 6).Value = 4280
$ws.Cells.Item(49, 6).Value = 170

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 18
$ws.Cells.Item(4, 6).Value = 15
$ws.Cells.Item(7, 6).Value = 2
$ws.Cells.Item(8, 6).Value = 34
$ws.Cells.Item(9, 6).Value = 66
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(18, 6).Value = 6
$ws.Cells.Item(19, 6).Value = 11
$ws.Cells.Item(21, 6).Value = 5
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(23, 6).Value = 74
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(30, 6).Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 6430

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 548
$ws.Cells.Item(3, 6).Value = 227
$ws.Cells.Item(4, 6).Value = 94
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 12042
$ws.Cells.Item(10, 6).Value = 12615
$ws.Cells.Item(12, 6).Value = 1299
$ws.Cells.Item(13, 6).Value = 1275
$ws.Cells.Item(17, 6).Value = 75
$ws.Cells.Item(18, 6).Value = 1418
$ws.Cells.Item(20, 6).Value = 1513
$ws.Cells.Item(21, 6).Value = 878
$ws.Cells.Item(22, 6).Value = 5
$ws.Cells.Item(23, 6).Value = 13
$ws.Cells.Item(25, 6).Value = 714
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 2003
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 12
$ws.Cells.Item(35, 6).Value = 74
$ws.Cells.Item(36, 6).Value = 33
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(40, 6).Value = 124
$ws.Cells.Item(43, 6).Value = 912
$ws.Cells.Item(44, 6).Value = 281
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(47, 6).Value = 170
$ws.Cells.Item(48, 6).Value = 0

